$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.939.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.43%  "
$ws.Range("D3").Value = "'1.866.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'318.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4375"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.84%  "
$ws.Range("D8").Value = "'0.3708"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").Value = "'0.07506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "'0.9373"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.37%  "
$ws.Range("D11").Value = "'21.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").Value = "'1.898.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").Value = "'6.741"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "'5.454"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").Value = "'0.06842"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").Value = "'0.000009071"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.12%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'15.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("D21").Value = "'27.930.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.39%  "
$ws.Range("D22").Value = "'5.108"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("D23").Value = "'11.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'2.122.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'2.004"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("D26").Value = "'153.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("D28").Value = "'5.495"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("D29").Value = "'113.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("D30").Value = "'1.722"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.29%  "
$ws.Range("D31").Value = "'0.09029"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").Value = "'0.8146"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.65%  "
$ws.Range("D33").Value = "'4.822"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.39%  "
$ws.Range("E34").Value = "  -5.70%  "
$ws.Range("D35").Value = "'2.956"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").Value = "'1.003"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "'0.05499"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.90%  "
$ws.Range("D38").Value = "'1.120"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").Value = "'0.01980"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.00%  "
$ws.Range("D40").Value = "'2.950"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "'0.5265"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("D42").Value = "'7.054"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.80%  "
$ws.Range("D43").Value = "'0.1704"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").Value = "'8.803"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.98%  "
$ws.Range("D45").Value = "'0.06773"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").Value = "'0.4915"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.89%  "
$ws.Range("D47").Value = "'10.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.42%  "
$ws.Range("D48").Value = "'107.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("E49").Value = "  -5.57%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "'1.883"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -13.19%  "
